$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "500 Global Flagship VC (non-accelerator checks)",
    "U.S. Dept. of Education (EIR Program)",
    "Owl Ventures",
    "Reach Capital",
    "Bill & Melinda Gates Foundation",
    "Google for Education",
    "National Science Foundation (NSF)",
    "500 Global (seed/accelerator)",
    "IES SBIR (ED/IES)",
    "Berkeley SkyDeck Fund (UC Berkeley)",
    "Y Combinator",
    "TGR Foundation (Tiger Woods)",
    "NSF SBIR (“America’s Seed Fund”)",
    "MIT Solve (Global Learning & education tracks)",
    "Washington Commanders Foundation",
    "Golden State Warriors Community Foundation",
    "NFL Foundation",
    "Indiana Pacers Foundation",
    "Chicago Blackhawks Foundation",
    "D.C. United Foundation",
    "Chicago Bulls Charities",
    "San Antonio Spurs - Spurs Give",
    "Nike Community Impact Fund",
    "Toronto FC - MLSE Foundation",
    "Houston Texans Foundation",
    "Seattle Sounders FC RAVE Foundation",
    "Minnesota Wild Foundation",
    "Las Vegas Raiders Foundation",
    "Detroit Red Wings Foundation",
    "Boston Red Sox Foundation"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
